$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column B values: "Gl_Nap" -> "[gl_nap1]"
$ws.Range("B2:B11").Value = "[gl_nap1]"

# 2) Insert a new column before column X so the old X:AB block shifts to Y:AC
$ws.Range("X1").EntireColumn.Insert()

# 3) New column header + values ("Vergleich" column, filled with "j")
$ws.Range("X1").Value = "Vergleich"
$ws.Range("X2:X11").Value = "j"

# 4) Update the view: scroll right and move the active selection
$ws.Range("V4").Select()
$excel.ActiveWindow.ScrollColumn = 13
